# Auto-generated edit script for cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.974.90"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "1.554.94"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.45%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.482"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.68"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.01%  "
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0589"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.25%  "
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "1.775.83"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "1.553.84"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.516"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").Value = "26.965.26"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "0.0₃0689"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("E23").Value = "  +1.72%  "
$ws.Range("E24").Value = "  -1.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.103"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0463"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("E32").Value = "  +1.21%  "
$ws.Range("D33").Value = "1.387.86"
$ws.Range("E33").Value = "  +2.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.03%  "
$ws.Range("E35").Value = "  +3.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.967"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.50%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.525"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.811"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.08%  "
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.990"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.25"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("D47").Value = "1.689.89"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0511"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0955"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.55%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.31%  "
